function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 5
Set-TextCell $ws "B2" "Hola te amo mucho"
Set-TextCell $ws "C2" "2024-09-08"
Set-TextCell $ws "D2" "2"
Set-TextCell $ws "E2" "5"
$ws.Range("F2").Value = 5
Set-TextCell $ws "G2" "M1B"
Set-TextCell $ws "H2" "Práctica"
Set-TextCell $ws "I2" "Ndjfjdndjfjd"
Set-TextCell $ws "J2" "Djfjdjtjdj"
Set-TextCell $ws "K2" "Jfjfjfjdd"

# Row 3
$ws.Range("A3").Value = 6
Set-TextCell $ws "B3" "aa"
Set-TextCell $ws "C3" "2024-09-13"
Set-TextCell $ws "D3" "fasdf"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 43434
Set-TextCell $ws "G3" "6a"
Set-TextCell $ws "H3" "sdgsdfg"
Set-TextCell $ws "I3" "dfg"
Set-TextCell $ws "J3" "sdfgsd"
Set-TextCell $ws "K3" "fgsdfgsdf"

# Row 4
$ws.Range("A4").Value = 7
Set-TextCell $ws "B4" "asdfasd"
Set-TextCell $ws "C4" "223223-02-23"
Set-TextCell $ws "D4" "as2"
$ws.Range("E4").Value = 23
$ws.Range("F4").Value = 232
Set-TextCell $ws "G4" "6a"
Set-TextCell $ws "H4" "sadfsa"
Set-TextCell $ws "I4" "dfasd"
Set-TextCell $ws "J4" "fasdf"
Set-TextCell $ws "K4" "asf"

# Row 5
$ws.Range("A5").Value = 9
Set-TextCell $ws "B5" "AAAAAAAAAAA"
Set-TextCell $ws "C5" "2024-09-08"
Set-TextCell $ws "D5" "2423423"
$ws.Range("E5").Value = 23423423
$ws.Range("F5").Value = 34243
Set-TextCell $ws "G5" "6a"
Set-TextCell $ws "H5" "ERWER"
Set-TextCell $ws "I5" "FDSFASDF"
Set-TextCell $ws "J5" "FDSFDS"
Set-TextCell $ws "K5" "FDSFSFDSF"

# Row 6
$ws.Range("A6").Value = 8
Set-TextCell $ws "B6" "23"
Set-TextCell $ws "C6" "2024-09-12"
Set-TextCell $ws "D6" ""
$ws.Range("E6").Value = 23
$ws.Range("F6").Value = 23
Set-TextCell $ws "G6" "M1B"
Set-TextCell $ws "H6" "ASDF"
Set-TextCell $ws "I6" "ASFASDF"
Set-TextCell $ws "J6" "ASDFAS"
Set-TextCell $ws "K6" "DFA"

# Row 7
$ws.Range("A7").Value = 9
Set-TextCell $ws "B7" "23"
Set-TextCell $ws "C7" "2024-09-12"
Set-TextCell $ws "D7" ""
$ws.Range("E7").Value = 232
$ws.Range("F7").Value = 232
Set-TextCell $ws "G7" "S6A"
Set-TextCell $ws "H7" "Práctica"
Set-TextCell $ws "I7" "sdfa"
Set-TextCell $ws "J7" "sdfas"
Set-TextCell $ws "K7" "dfas"

# Row 8
$ws.Range("A8").Value = 10
Set-TextCell $ws "B8" "asdf"
Set-TextCell $ws "C8" "2024-09-12"
Set-TextCell $ws "D8" ""
Set-TextCell $ws "E8" "23"
Set-TextCell $ws "F8" ""
Set-TextCell $ws "G8" "M1B"
Set-TextCell $ws "H8" "Teórico/Práctica"
Set-TextCell $ws "I8" "asdf"
Set-TextCell $ws "J8" "asdfsa"
Set-TextCell $ws "K8" "dfsaddf"

# Row 9
$ws.Range("A9").Value = 11
Set-TextCell $ws "B9" "234"
Set-TextCell $ws "C9" "2024-09-12"
Set-TextCell $ws "D9" ""
Set-TextCell $ws "E9" "23"
Set-TextCell $ws "F9" ""
Set-TextCell $ws "G9" "TM"
Set-TextCell $ws "H9" "Explicación"
Set-TextCell $ws "I9" "safd"
Set-TextCell $ws "J9" "asdf"
Set-TextCell $ws "K9" "asdfas"
